$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Navigation")

# Rename the two section headers that lost the "User " prefix
$ws.Range("A2").Value = "Site Menu Navigation"
$ws.Range("A11").Value = "Site Sub Menu Transmittals Navigation"

# Insert a new row for the "Home" menu link right after the
# "Document & File Storage" entry (old row 9), before the
# "Transmittals" sub-menu section header (old row 11).
$ws.Rows("10:10").Insert()

# The inserted row copies formatting from the row above (row 9, which
# uses the red "last item" style) - reset it back to the default style
# so the new cells carry no explicit style, matching a plain data row.
$ws.Range("B10:E10").Style = "Normal"

$ws.Range("B10").Value = "Menu - Home"
$ws.Range("C10").Value = "xpath"
$ws.Range("D10").Value = "link"
$ws.Range("E10").Value = "Home"

# Match the saved selection state (active cell on the new row).
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
